$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert 13 new rows after row 15 (before the "* Footer" row), pushing the
# footer from row 16 down to row 29. The newly inserted rows inherit the
# formatting of row 15 (the row above), which already matches the blank
# "spacer" style used by the new row 16 and the new row 22 / 28 spacer rows.
# ---------------------------------------------------------------------------
$ws.Range("A16:A28").EntireRow.Insert()

# Row 16 is a blank spacer row (same look as rows 7/15) - nothing else to do,
# it already inherited the correct formatting from the insert above.

# Clear out the column-A gray highlight on rows that should not have it
# (only A19 and A25 keep the directive-marker look in this new block).
$ws.Range("A16:A18").ClearContents()
$ws.Range("A16:A18").Interior.Pattern = -4142
$ws.Range("A20:A24").ClearContents()
$ws.Range("A20:A24").Interior.Pattern = -4142
$ws.Range("A26:A28").ClearContents()
$ws.Range("A26:A28").Interior.Pattern = -4142

# ---------------------------------------------------------------------------
# Helper block writer: builds the "m1..m5 / m1-2 .. m4-5" merged-cell demo
# grid used to test that merged cells survive row insertion/removal.
# Block 1 occupies rows 17-21 (+22 spacer); block 2 occupies rows 23-27 (+28).
# ---------------------------------------------------------------------------
function New-MergeDemoBlock($r0, $markerText) {
    $r1 = $r0
    $r2 = $r0 + 1
    $r3 = $r0 + 2
    $r4 = $r0 + 3
    $r5 = $r0 + 4

    $ws.Range("B$r1").Value = "m1"
    $ws.Range("B$r2").Value = "m2"
    $ws.Range("B$r3").Value = "m3"
    $ws.Range("B$r4").Value = "m4"
    $ws.Range("B$r5").Value = "m5"

    $ws.Range("A$r3").Value = $markerText

    $ws.Range("E$r1`:E$r2").Merge()
    $ws.Range("E$r1`:E$r2").Value = "m1-2"

    $ws.Range("F$r2`:F$r3").Merge()
    $ws.Range("F$r2`:F$r3").Value = "m2-3"

    $ws.Range("G$r2`:G$r4").Merge()
    $ws.Range("G$r2`:G$r4").Value = "m2-4"

    $ws.Range("C$r3`:D$r3").Merge()
    $ws.Range("C$r3`:D$r3").Value = "m3"

    $ws.Range("H$r3`:H$r4").Merge()
    $ws.Range("H$r3`:H$r4").Value = "m3-4"

    $ws.Range("I$r4`:I$r5").Merge()
    $ws.Range("I$r4`:I$r5").Value = "m4-5"

    # Outer box around the whole grid, then a box around every merged range
    # so that the internal divider lines line up correctly.
    $ws.Range("B$r1`:I$r5").BorderAround(1)
    $ws.Range("E$r1`:E$r2").BorderAround(1)
    $ws.Range("F$r2`:F$r3").BorderAround(1)
    $ws.Range("G$r2`:G$r4").BorderAround(1)
    $ws.Range("C$r3`:D$r3").BorderAround(1)
    $ws.Range("H$r3`:H$r4").BorderAround(1)
    $ws.Range("I$r4`:I$r5").BorderAround(1)

    # Merged "label" cells are centered both ways.
    $ws.Range("E$r1`:E$r2").HorizontalAlignment = -4108
    $ws.Range("E$r1`:E$r2").VerticalAlignment = -4108
    $ws.Range("F$r2`:F$r3").HorizontalAlignment = -4108
    $ws.Range("F$r2`:F$r3").VerticalAlignment = -4108
    $ws.Range("G$r2`:G$r4").HorizontalAlignment = -4108
    $ws.Range("G$r2`:G$r4").VerticalAlignment = -4108
    $ws.Range("C$r3`:D$r3").HorizontalAlignment = -4108
    $ws.Range("C$r3`:D$r3").VerticalAlignment = -4108
    $ws.Range("H$r3`:H$r4").HorizontalAlignment = -4108
    $ws.Range("H$r3`:H$r4").VerticalAlignment = -4108
    $ws.Range("I$r4`:I$r5").HorizontalAlignment = -4108
    $ws.Range("I$r4`:I$r5").VerticalAlignment = -4108
}

New-MergeDemoBlock 17 "{{#if awards}}"
New-MergeDemoBlock 23 "{{#each works}}"

# ---------------------------------------------------------------------------
# Misc cosmetic / view metadata changes captured by the diff.
# ---------------------------------------------------------------------------
$ws.Range("A25").Select()
$excel.ActiveWindow.WindowState = $excel.ActiveWindow.WindowState
